$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently holds weekly "Kiwi / Hayward" price records for
# "Feria Lagunitas de Puerto Montt" in rows 437-452 (A1:T452).
# A new week of data (3 rows: Especial/Primera/Segunda) needs to be
# inserted at the top of that block, pushing the existing 16 rows down
# by 3 (new dimension A1:T455).

$ws.Range("A437:T439").EntireRow.Insert()

# Row 437 - Especial
$ws.Range("A437").Value = 4
$ws.Range("B437").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C437").Value = "Los Lagos"
$ws.Range("D437").Value = 45008
$ws.Range("E437").Value = 10
$ws.Range("F437").Value = "Fruta"
$ws.Range("G437").Value = 100101
$ws.Range("H437").Value = "Berries"
$ws.Range("I437").Value = 100101007
$ws.Range("J437").Value = "Kiwi"
$ws.Range("K437").Value = "Hayward"
$ws.Range("L437").Value = "Especial"
$ws.Range("M437").Value = 200
$ws.Range("N437").Value = 22000
$ws.Range("O437").Value = 22000
$ws.Range("P437").Value = 22000
$ws.Range("Q437").Value = "$/caja 15 kilos"
$ws.Range("R437").Value = "Región de O'Higgins"
$ws.Range("S437").Value = 1467
$ws.Range("T437").Value = 15

# Row 438 - Primera
$ws.Range("A438").Value = 4
$ws.Range("B438").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C438").Value = "Los Lagos"
$ws.Range("D438").Value = 45008
$ws.Range("E438").Value = 10
$ws.Range("F438").Value = "Fruta"
$ws.Range("G438").Value = 100101
$ws.Range("H438").Value = "Berries"
$ws.Range("I438").Value = 100101007
$ws.Range("J438").Value = "Kiwi"
$ws.Range("K438").Value = "Hayward"
$ws.Range("L438").Value = "Primera"
$ws.Range("M438").Value = 200
$ws.Range("N438").Value = 19000
$ws.Range("O438").Value = 19000
$ws.Range("P438").Value = 19000
$ws.Range("Q438").Value = "$/caja 15 kilos"
$ws.Range("R438").Value = "Región de O'Higgins"
$ws.Range("S438").Value = 1267
$ws.Range("T438").Value = 15

# Row 439 - Segunda
$ws.Range("A439").Value = 4
$ws.Range("B439").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C439").Value = "Los Lagos"
$ws.Range("D439").Value = 45008
$ws.Range("E439").Value = 10
$ws.Range("F439").Value = "Fruta"
$ws.Range("G439").Value = 100101
$ws.Range("H439").Value = "Berries"
$ws.Range("I439").Value = 100101007
$ws.Range("J439").Value = "Kiwi"
$ws.Range("K439").Value = "Hayward"
$ws.Range("L439").Value = "Segunda"
$ws.Range("M439").Value = 200
$ws.Range("N439").Value = 16000
$ws.Range("O439").Value = 16000
$ws.Range("P439").Value = 16000
$ws.Range("Q439").Value = "$/caja 15 kilos"
$ws.Range("R439").Value = "Región de O'Higgins"
$ws.Range("S439").Value = 1067
$ws.Range("T439").Value = 15

Write-Output ("Dimension: " + $ws.UsedRange.Address())
